$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-5 with the corrected data
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "Supermercados Rey"
$ws.Range("C2").Value = 1

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "Supermercados Rey"
$ws.Range("C3").Value = 3

$ws.Range("A4").Value = 6
$ws.Range("B4").Value = "Embonor"
$ws.Range("C4").Value = 6

$ws.Range("A5").Value = 8
$ws.Range("B5").Value = "Embonor"
$ws.Range("C5").Value = 3

# Remove the now-obsolete row 6
$ws.Rows.Item(6).Delete()
